$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3368
$ws1.Range("F5").Value = 6950
$ws1.Range("F6").Value = 2365
$ws1.Range("F8").Value = 100

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3368
$ws4.Range("F6").Value = 6950
$ws4.Range("F7").Value = 2365
$ws4.Range("F9").Value = 100
